$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append 4 new survey response rows (131-134) ---
# Row 131
$ws.Cells.Item(131,1).Value = 25
$ws.Cells.Item(131,2).Value = 1
$ws.Cells.Item(131,3).Value = 6
$ws.Cells.Item(131,4).Value = 6
$ws.Cells.Item(131,5).Value = 1
$ws.Cells.Item(131,6).Value = 1
$ws.Cells.Item(131,7).Value = 2
$ws.Cells.Item(131,8).Value = 0
$ws.Cells.Item(131,10).Value = "Aschaffenburg"
$ws.Cells.Item(131,11).Value = 0
$ws.Cells.Item(131,12).Value = 0
$ws.Cells.Item(131,13).Value = 0
$ws.Cells.Item(131,14).Value = 0
$ws.Cells.Item(131,15).Value = 6
$ws.Cells.Item(131,16).Value = 4
$ws.Cells.Item(131,17).Value = 3
$ws.Cells.Item(131,18).Value = 4
$ws.Cells.Item(131,19).Value = 6
$ws.Cells.Item(131,20).Value = 3
$ws.Cells.Item(131,21).Value = 4
$ws.Cells.Item(131,22).Value = 5
$ws.Cells.Item(131,23).Value = 7
$ws.Cells.Item(131,24).Value = 4
$ws.Cells.Item(131,25).Value = 4
$ws.Cells.Item(131,26).Value = 4
$ws.Cells.Item(131,27).Value = 5
$ws.Cells.Item(131,28).Value = 6
$ws.Cells.Item(131,29).Value = 4
$ws.Cells.Item(131,30).Value = 5
$ws.Cells.Item(131,31).Value = 4
$ws.Cells.Item(131,32).Value = 7
$ws.Cells.Item(131,33).Value = 3
$ws.Cells.Item(131,34).Value = 6
$ws.Cells.Item(131,35).Value = 6
$ws.Cells.Item(131,36).Value = 6
$ws.Cells.Item(131,37).Value = 6
$ws.Cells.Item(131,38).Value = 6
$ws.Cells.Item(131,39).Value = 6
$ws.Cells.Item(131,40).Value = 7
$ws.Cells.Item(131,41).Value = 7
$ws.Cells.Item(131,42).Value = 7
$ws.Cells.Item(131,43).Value = 6
$ws.Cells.Item(131,44).Value = 5
$ws.Cells.Item(131,45).Value = 6
$ws.Cells.Item(131,46).Value = 6
$ws.Cells.Item(131,47).Value = 6
$ws.Cells.Item(131,48).Value = 6
$ws.Cells.Item(131,49).Value = 4
$ws.Cells.Item(131,50).Value = 6
$ws.Cells.Item(131,51).Value = 5
$ws.Cells.Item(131,52).Value = 7
$ws.Cells.Item(131,53).Value = 5
$ws.Cells.Item(131,54).Value = 6
$ws.Cells.Item(131,55).Value = 5
$ws.Cells.Item(131,56).Value = 5
$ws.Cells.Item(131,57).Value = 3
$ws.Cells.Item(131,58).Value = 5
$ws.Cells.Item(131,59).Value = 5
$ws.Cells.Item(131,60).Value = 6
$ws.Cells.Item(131,61).Value = 6
$ws.Cells.Item(131,62).Value = 3
$ws.Cells.Item(131,63).Value = 5
$ws.Cells.Item(131,64).Value = 1
$ws.Cells.Item(131,65).Value = 3
$ws.Cells.Item(131,66).Value = 5
$ws.Cells.Item(131,67).Value = 5
$ws.Cells.Item(131,68).Value = 4
$ws.Cells.Item(131,69).Value = 7
$ws.Cells.Item(131,70).Value = 7
$ws.Cells.Item(131,71).Value = 7
$ws.Cells.Item(131,72).Value = 4
$ws.Cells.Item(131,73).Value = 4
$ws.Cells.Item(131,74).Value = 4
$ws.Cells.Item(131,75).Value = 4
$ws.Cells.Item(131,76).Value = 4
$ws.Cells.Item(131,77).Value = 4
$ws.Cells.Item(131,78).Value = 17.438666666666666

# Row 132
$ws.Cells.Item(132,1).Value = 33
$ws.Cells.Item(132,2).Value = 2
$ws.Cells.Item(132,3).Value = 6
$ws.Cells.Item(132,4).Value = 5
$ws.Cells.Item(132,5).Value = 1
$ws.Cells.Item(132,6).Value = 1
$ws.Cells.Item(132,7).Value = 4
$ws.Cells.Item(132,8).Value = 1
$ws.Cells.Item(132,9).Value = 6
$ws.Cells.Item(132,10).Value = "Wietzendorf"
$ws.Cells.Item(132,11).Value = 2
$ws.Cells.Item(132,12).Value = 0
$ws.Cells.Item(132,13).Value = 0
$ws.Cells.Item(132,14).Value = 0
$ws.Cells.Item(132,15).Value = 6
$ws.Cells.Item(132,16).Value = 2
$ws.Cells.Item(132,17).Value = 5
$ws.Cells.Item(132,18).Value = 4
$ws.Cells.Item(132,19).Value = 6
$ws.Cells.Item(132,20).Value = 5
$ws.Cells.Item(132,21).Value = 5
$ws.Cells.Item(132,22).Value = 6
$ws.Cells.Item(132,23).Value = 7
$ws.Cells.Item(132,24).Value = 4
$ws.Cells.Item(132,25).Value = 3
$ws.Cells.Item(132,26).Value = 6
$ws.Cells.Item(132,27).Value = 4
$ws.Cells.Item(132,28).Value = 6
$ws.Cells.Item(132,29).Value = 5
$ws.Cells.Item(132,30).Value = 5
$ws.Cells.Item(132,31).Value = 6
$ws.Cells.Item(132,32).Value = 7
$ws.Cells.Item(132,33).Value = 6
$ws.Cells.Item(132,34).Value = 5
$ws.Cells.Item(132,35).Value = 6
$ws.Cells.Item(132,36).Value = 4
$ws.Cells.Item(132,37).Value = 3
$ws.Cells.Item(132,38).Value = 2
$ws.Cells.Item(132,39).Value = 5
$ws.Cells.Item(132,40).Value = 6
$ws.Cells.Item(132,41).Value = 5
$ws.Cells.Item(132,42).Value = 5
$ws.Cells.Item(132,43).Value = 5
$ws.Cells.Item(132,44).Value = 6
$ws.Cells.Item(132,45).Value = 5
$ws.Cells.Item(132,46).Value = 6
$ws.Cells.Item(132,47).Value = 6
$ws.Cells.Item(132,48).Value = 6
$ws.Cells.Item(132,49).Value = 4
$ws.Cells.Item(132,50).Value = 2
$ws.Cells.Item(132,51).Value = 3
$ws.Cells.Item(132,52).Value = 5
$ws.Cells.Item(132,53).Value = 3
$ws.Cells.Item(132,54).Value = 4
$ws.Cells.Item(132,55).Value = 5
$ws.Cells.Item(132,56).Value = 4
$ws.Cells.Item(132,57).Value = 2
$ws.Cells.Item(132,58).Value = 5
$ws.Cells.Item(132,59).Value = 3
$ws.Cells.Item(132,60).Value = 3
$ws.Cells.Item(132,61).Value = 1
$ws.Cells.Item(132,62).Value = 5
$ws.Cells.Item(132,63).Value = 4
$ws.Cells.Item(132,64).Value = 1
$ws.Cells.Item(132,65).Value = 2
$ws.Cells.Item(132,66).Value = 5
$ws.Cells.Item(132,67).Value = 5
$ws.Cells.Item(132,68).Value = 5
$ws.Cells.Item(132,69).Value = 6
$ws.Cells.Item(132,70).Value = 6
$ws.Cells.Item(132,71).Value = 6
$ws.Cells.Item(132,72).Value = 2
$ws.Cells.Item(132,73).Value = 2
$ws.Cells.Item(132,74).Value = 3
$ws.Cells.Item(132,75).Value = 2
$ws.Cells.Item(132,76).Value = 4
$ws.Cells.Item(132,77).Value = 2
$ws.Cells.Item(132,78).Value = 18.89

# Row 133
$ws.Cells.Item(133,1).Value = 35
$ws.Cells.Item(133,2).Value = 2
$ws.Cells.Item(133,3).Value = 6
$ws.Cells.Item(133,4).Value = 5
$ws.Cells.Item(133,5).Value = 1
$ws.Cells.Item(133,6).Value = 1
$ws.Cells.Item(133,7).Value = 2
$ws.Cells.Item(133,8).Value = 1
$ws.Cells.Item(133,9).Value = 5
$ws.Cells.Item(133,10).Value = "Konstanz"
$ws.Cells.Item(133,11).Value = 1
$ws.Cells.Item(133,12).Value = 1
$ws.Cells.Item(133,13).Value = 1
$ws.Cells.Item(133,14).Value = 0
$ws.Cells.Item(133,15).Value = 6
$ws.Cells.Item(133,16).Value = 5
$ws.Cells.Item(133,17).Value = 4
$ws.Cells.Item(133,18).Value = 6
$ws.Cells.Item(133,19).Value = 7
$ws.Cells.Item(133,20).Value = 3
$ws.Cells.Item(133,21).Value = 5
$ws.Cells.Item(133,22).Value = 6
$ws.Cells.Item(133,23).Value = 7
$ws.Cells.Item(133,24).Value = 4
$ws.Cells.Item(133,25).Value = 5
$ws.Cells.Item(133,26).Value = 4
$ws.Cells.Item(133,27).Value = 6
$ws.Cells.Item(133,28).Value = 6
$ws.Cells.Item(133,29).Value = 3
$ws.Cells.Item(133,30).Value = 5
$ws.Cells.Item(133,31).Value = 5
$ws.Cells.Item(133,32).Value = 7
$ws.Cells.Item(133,33).Value = 4
$ws.Cells.Item(133,34).Value = 6
$ws.Cells.Item(133,35).Value = 3
$ws.Cells.Item(133,36).Value = 4
$ws.Cells.Item(133,37).Value = 6
$ws.Cells.Item(133,38).Value = 6
$ws.Cells.Item(133,39).Value = 7
$ws.Cells.Item(133,40).Value = 7
$ws.Cells.Item(133,41).Value = 5
$ws.Cells.Item(133,42).Value = 6
$ws.Cells.Item(133,43).Value = 6
$ws.Cells.Item(133,44).Value = 4
$ws.Cells.Item(133,45).Value = 7
$ws.Cells.Item(133,46).Value = 6
$ws.Cells.Item(133,47).Value = 6
$ws.Cells.Item(133,48).Value = 6
$ws.Cells.Item(133,49).Value = 4
$ws.Cells.Item(133,50).Value = 6
$ws.Cells.Item(133,51).Value = 5
$ws.Cells.Item(133,52).Value = 7
$ws.Cells.Item(133,53).Value = 4
$ws.Cells.Item(133,54).Value = 6
$ws.Cells.Item(133,55).Value = 6
$ws.Cells.Item(133,56).Value = 5
$ws.Cells.Item(133,57).Value = 3
$ws.Cells.Item(133,58).Value = 4
$ws.Cells.Item(133,59).Value = 2
$ws.Cells.Item(133,60).Value = 5
$ws.Cells.Item(133,61).Value = 4
$ws.Cells.Item(133,62).Value = 3
$ws.Cells.Item(133,63).Value = 5
$ws.Cells.Item(133,64).Value = 2
$ws.Cells.Item(133,65).Value = 1
$ws.Cells.Item(133,66).Value = 3
$ws.Cells.Item(133,67).Value = 6
$ws.Cells.Item(133,68).Value = 6
$ws.Cells.Item(133,69).Value = 6
$ws.Cells.Item(133,70).Value = 6
$ws.Cells.Item(133,71).Value = 6
$ws.Cells.Item(133,72).Value = 3
$ws.Cells.Item(133,73).Value = 3
$ws.Cells.Item(133,74).Value = 4
$ws.Cells.Item(133,75).Value = 3
$ws.Cells.Item(133,76).Value = 4
$ws.Cells.Item(133,77).Value = 5
$ws.Cells.Item(133,78).Value = 32.15

# Row 134
$ws.Cells.Item(134,1).Value = 37
$ws.Cells.Item(134,2).Value = 2
$ws.Cells.Item(134,3).Value = 6
$ws.Cells.Item(134,4).Value = 10
$ws.Cells.Item(134,5).Value = 1
$ws.Cells.Item(134,6).Value = 3
$ws.Cells.Item(134,7).Value = 8
$ws.Cells.Item(134,8).Value = 1
$ws.Cells.Item(134,9).Value = 1
$ws.Cells.Item(134,10).Value = "Leipzig"
$ws.Cells.Item(134,11).Value = 7
$ws.Cells.Item(134,12).Value = 1
$ws.Cells.Item(134,13).Value = 1
$ws.Cells.Item(134,14).Value = 0
$ws.Cells.Item(134,15).Value = 6
$ws.Cells.Item(134,16).Value = 3
$ws.Cells.Item(134,17).Value = 3
$ws.Cells.Item(134,18).Value = 4
$ws.Cells.Item(134,19).Value = 6
$ws.Cells.Item(134,20).Value = 3
$ws.Cells.Item(134,21).Value = 4
$ws.Cells.Item(134,22).Value = 4
$ws.Cells.Item(134,23).Value = 6
$ws.Cells.Item(134,24).Value = 4
$ws.Cells.Item(134,25).Value = 2
$ws.Cells.Item(134,26).Value = 3
$ws.Cells.Item(134,27).Value = 4
$ws.Cells.Item(134,28).Value = 5
$ws.Cells.Item(134,29).Value = 3
$ws.Cells.Item(134,30).Value = 4
$ws.Cells.Item(134,31).Value = 4
$ws.Cells.Item(134,32).Value = 6
$ws.Cells.Item(134,33).Value = 5
$ws.Cells.Item(134,34).Value = 2
$ws.Cells.Item(134,35).Value = 2
$ws.Cells.Item(134,36).Value = 4
$ws.Cells.Item(134,37).Value = 3
$ws.Cells.Item(134,38).Value = 3
$ws.Cells.Item(134,39).Value = 5
$ws.Cells.Item(134,40).Value = 6
$ws.Cells.Item(134,41).Value = 7
$ws.Cells.Item(134,42).Value = 6
$ws.Cells.Item(134,43).Value = 4
$ws.Cells.Item(134,44).Value = 6
$ws.Cells.Item(134,45).Value = 7
$ws.Cells.Item(134,46).Value = 7
$ws.Cells.Item(134,47).Value = 7
$ws.Cells.Item(134,48).Value = 6
$ws.Cells.Item(134,49).Value = 4
$ws.Cells.Item(134,50).Value = 5
$ws.Cells.Item(134,51).Value = 5
$ws.Cells.Item(134,52).Value = 5
$ws.Cells.Item(134,53).Value = 6
$ws.Cells.Item(134,54).Value = 5
$ws.Cells.Item(134,55).Value = 4
$ws.Cells.Item(134,56).Value = 4
$ws.Cells.Item(134,57).Value = 4
$ws.Cells.Item(134,58).Value = 5
$ws.Cells.Item(134,59).Value = 5
$ws.Cells.Item(134,60).Value = 5
$ws.Cells.Item(134,61).Value = 4
$ws.Cells.Item(134,62).Value = 3
$ws.Cells.Item(134,63).Value = 6
$ws.Cells.Item(134,64).Value = 2
$ws.Cells.Item(134,65).Value = 2
$ws.Cells.Item(134,66).Value = 4
$ws.Cells.Item(134,67).Value = 3
$ws.Cells.Item(134,68).Value = 3
$ws.Cells.Item(134,69).Value = 6
$ws.Cells.Item(134,70).Value = 6
$ws.Cells.Item(134,71).Value = 6
$ws.Cells.Item(134,72).Value = 5
$ws.Cells.Item(134,73).Value = 5
$ws.Cells.Item(134,74).Value = 5
$ws.Cells.Item(134,75).Value = 5
$ws.Cells.Item(134,76).Value = 5
$ws.Cells.Item(134,77).Value = 4
$ws.Cells.Item(134,78).Value = 14.768666666666666

# --- Adjust column AA (27) width to match updated best-fit content ---
$ws.Columns.Item(27).ColumnWidth = 14

# --- Restore the last on-screen selection ---
$ws.Range("Z141").Select()
